# Master_Gantt.xlsx update: "Updated hours, added comments"
#
# Devin D.'s (column C) task list was reworked (the "Program
# core/optional features / Debug main/optional features / Documentation /
# Testing / Release Build" placeholders were replaced with the real task
# breakdown for the obstacles/level-exit work), and her actual logged
# "Duration" / "Completed" hours for each of the 9 tasks were corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gantt")

# --- Task names (column C, rows 3-10) -------------------------------------
# NB: cells are written in the same order the new unique strings were
# first introduced in the workbook (C5,C6,C7,C8,C9,C3,C10) so the
# shared-string table is rebuilt in the same sequence; C4 re-uses a string
# that already exists ("Write Public functions, basic structure").
$ws.Range("C5").Value  = "Make inert obstacles"
$ws.Range("C6").Value  = "Make Simple Interactables"
$ws.Range("C7").Value  = "Make Level Exit/Dodo Egg"
$ws.Range("C8").Value  = "Make Optional Obstacles"
$ws.Range("C9").Value  = "Integrate other's parts that were not previously working"
$ws.Range("C3").Value  = "Github setup, Etc."
$ws.Range("C10").Value = "Final testing and build"
$ws.Range("C4").Value  = "Write Public functions, basic structure"

# --- Duration hours (column C, rows 20-27; row 28 unchanged) --------------
$ws.Range("C20").Value = 24
$ws.Range("C21").Value = 6
$ws.Range("C22").Value = 6
$ws.Range("C23").Value = 6
$ws.Range("C24").Value = 3
$ws.Range("C25").Value = 3
$ws.Range("C26").Value = 6
$ws.Range("C27").Value = 3

# --- Completed hours (column C, rows 29-37) --------------------------------
$ws.Range("C29").Value = 24
$ws.Range("C30").Value = 18
$ws.Range("C31").Value = 7
$ws.Range("C32").Value = 4
$ws.Range("C33").Value = 4
$ws.Range("C34").Value = 6
$ws.Range("C35").Value = 4
$ws.Range("C36").Value = 7
$ws.Range("C37").Value = 3

# --- View state: zoom level and active selection ---------------------------
$excel.ActiveWindow.Zoom = 115
$ws.Range("D40").Select()
